{"js": "// Replace each \"A\u00d7B=C\" multiplication-table answer in the document with its\n// new value. The mapping below is built from the unified diff: every old\n// answer string is unique in the document, so a literal (non-wildcard) body\n// search-and-replace is unambiguous and order-independent.\nconst replacements = [\n  [\"153\u00d73=459\", \"263\u00d75=1315\"],\n  [\"701\u00d75=3505\", \"461\u00d72=922\"],\n  [\"499\u00d74=1996\", \"767\u00d75=3835\"],\n  [\"120\u00d79=1080\", \"751\u00d75=3755\"],\n  [\"534\u00d79=4806\", \"794\u00d74=3176\"],\n  [\"493\u00d72=986\", \"858\u00d76=5148\"],\n  [\"271\u00d75=1355\", \"487\u00d79=4383\"],\n  [\"942\u00d73=2826\", \"562\u00d76=3372\"],\n  [\"659\u00d74=2636\", \"804\u00d78=6432\"],\n  [\"424\u00d78=3392\", \"415\u00d74=1660\"],\n  [\"502\u00d73=1506\", \"304\u00d77=2128\"],\n  [\"408\u00d78=3264\", \"836\u00d75=4180\"],\n  [\"510\u00d78=4080\", \"460\u00d74=1840\"],\n  [\"785\u00d75=3925\", \"857\u00d79=7713\"],\n  [\"127\u00d78=1016\", \"651\u00d79=5859\"],\n  [\"310\u00d79=2790\", \"240\u00d74=960\"],\n  [\"608\u00d77=4256\", \"649\u00d79=5841\"],\n  [\"188\u00d77=1316\", \"390\u00d74=1560\"],\n  [\"600\u00d74=2400\", \"112\u00d74=448\"],\n  [\"930\u00d76=5580\", \"460\u00d72=920\"],\n  [\"878\u00d79=7902\", \"538\u00d76=3228\"],\n  [\"524\u00d76=3144\", \"831\u00d79=7479\"],\n  [\"584\u00d72=1168\", \"176\u00d73=528\"],\n  [\"648\u00d73=1944\", \"476\u00d79=4284\"],\n  [\"261\u00d77=1827\", \"751\u00d73=2253\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find expected text \"${oldText}\" in document.`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each \"A\u00d7B=C\" multiplication-table answer in the document with its\n# new value. The mapping below is built from the unified diff: every old\n# answer string is unique in the document, so a literal (non-wildcard)\n# Find/Replace over the whole document body is unambiguous and\n# order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"153\u00d73=459\";   New = \"263\u00d75=1315\" },\n    @{ Old = \"701\u00d75=3505\";  New = \"461\u00d72=922\" },\n    @{ Old = \"499\u00d74=1996\";  New = \"767\u00d75=3835\" },\n    @{ Old = \"120\u00d79=1080\";  New = \"751\u00d75=3755\" },\n    @{ Old = \"534\u00d79=4806\";  New = \"794\u00d74=3176\" },\n    @{ Old = \"493\u00d72=986\";   New = \"858\u00d76=5148\" },\n    @{ Old = \"271\u00d75=1355\";  New = \"487\u00d79=4383\" },\n    @{ Old = \"942\u00d73=2826\";  New = \"562\u00d76=3372\" },\n    @{ Old = \"659\u00d74=2636\";  New = \"804\u00d78=6432\" },\n    @{ Old = \"424\u00d78=3392\";  New = \"415\u00d74=1660\" },\n    @{ Old = \"502\u00d73=1506\";  New = \"304\u00d77=2128\" },\n    @{ Old = \"408\u00d78=3264\";  New = \"836\u00d75=4180\" },\n    @{ Old = \"510\u00d78=4080\";  New = \"460\u00d74=1840\" },\n    @{ Old = \"785\u00d75=3925\";  New = \"857\u00d79=7713\" },\n    @{ Old = \"127\u00d78=1016\";  New = \"651\u00d79=5859\" },\n    @{ Old = \"310\u00d79=2790\";  New = \"240\u00d74=960\" },\n    @{ Old = \"608\u00d77=4256\";  New = \"649\u00d79=5841\" },\n    @{ Old = \"188\u00d77=1316\";  New = \"390\u00d74=1560\" },\n    @{ Old = \"600\u00d74=2400\";  New = \"112\u00d74=448\" },\n    @{ Old = \"930\u00d76=5580\";  New = \"460\u00d72=920\" },\n    @{ Old = \"878\u00d79=7902\";  New = \"538\u00d76=3228\" },\n    @{ Old = \"524\u00d76=3144\";  New = \"831\u00d79=7479\" },\n    @{ Old = \"584\u00d72=1168\";  New = \"176\u00d73=528\" },\n    @{ Old = \"648\u00d73=1944\";  New = \"476\u00d79=4284\" },\n    @{ Old = \"261\u00d77=1827\";  New = \"751\u00d73=2253\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $true,       # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Could not find expected text '$($pair.Old)' in document.\"\n    }\n}\n"}
